# Apply replacements for each math expression cell in the table.
# Each (old, new) pair is unique and unambiguous in the document, so we
# use Find/Execute with Replace:=wdReplaceAll (2) to swap the text of
# the w:t run in-place, cell by cell, in document order.
$d = $word.ActiveDocument

$d.Content.Find.Execute("80-24=", $true, $true, $false, $false, $false, $true, 1, $false, "28+70=", 2) | Out-Null
$d.Content.Find.Execute("41+5=", $true, $true, $false, $false, $false, $true, 1, $false, "39+56=", 2) | Out-Null
$d.Content.Find.Execute("61+2=", $true, $true, $false, $false, $false, $true, 1, $false, "3+0=", 2) | Out-Null
$d.Content.Find.Execute("84-38=", $true, $true, $false, $false, $false, $true, 1, $false, "4+77=", 2) | Out-Null
$d.Content.Find.Execute("20+40=", $true, $true, $false, $false, $false, $true, 1, $false, "18+32=", 2) | Out-Null
$d.Content.Find.Execute("41+7=", $true, $true, $false, $false, $false, $true, 1, $false, "12-9=", 2) | Out-Null
$d.Content.Find.Execute("75-21=", $true, $true, $false, $false, $false, $true, 1, $false, "5+66=", 2) | Out-Null
$d.Content.Find.Execute("6+32=", $true, $true, $false, $false, $false, $true, 1, $false, "89-13=", 2) | Out-Null
$d.Content.Find.Execute("73-45=", $true, $true, $false, $false, $false, $true, 1, $false, "92-65=", 2) | Out-Null
$d.Content.Find.Execute("75-50=", $true, $true, $false, $false, $false, $true, 1, $false, "83+4=", 2) | Out-Null
$d.Content.Find.Execute("43-23=", $true, $true, $false, $false, $false, $true, 1, $false, "26+16=", 2) | Out-Null
$d.Content.Find.Execute("16+63=", $true, $true, $false, $false, $false, $true, 1, $false, "83-12=", 2) | Out-Null
$d.Content.Find.Execute("34+35=", $true, $true, $false, $false, $false, $true, 1, $false, "46+50=", 2) | Out-Null
$d.Content.Find.Execute("87-59=", $true, $true, $false, $false, $false, $true, 1, $false, "42+10=", 2) | Out-Null
$d.Content.Find.Execute("76-25=", $true, $true, $false, $false, $false, $true, 1, $false, "38+22=", 2) | Out-Null
$d.Content.Find.Execute("70-35=", $true, $true, $false, $false, $false, $true, 1, $false, "83+6=", 2) | Out-Null
$d.Content.Find.Execute("73-44=", $true, $true, $false, $false, $false, $true, 1, $false, "38+0=", 2) | Out-Null
$d.Content.Find.Execute("70-67=", $true, $true, $false, $false, $false, $true, 1, $false, "38+40=", 2) | Out-Null
$d.Content.Find.Execute("5+85=", $true, $true, $false, $false, $false, $true, 1, $false, "32+0=", 2) | Out-Null
$d.Content.Find.Execute("35+44=", $true, $true, $false, $false, $false, $true, 1, $false, "68-42=", 2) | Out-Null
$d.Content.Find.Execute("88-79=", $true, $true, $false, $false, $false, $true, 1, $false, "99-55=", 2) | Out-Null
$d.Content.Find.Execute("47+46=", $true, $true, $false, $false, $false, $true, 1, $false, "64+5=", 2) | Out-Null
$d.Content.Find.Execute("96-31=", $true, $true, $false, $false, $false, $true, 1, $false, "71-40=", 2) | Out-Null
$d.Content.Find.Execute("69-48=", $true, $true, $false, $false, $false, $true, 1, $false, "34-26=", 2) | Out-Null
$d.Content.Find.Execute("46-10=", $true, $true, $false, $false, $false, $true, 1, $false, "8-6=", 2) | Out-Null
$d.Content.Find.Execute("88-32=", $true, $true, $false, $false, $false, $true, 1, $false, "44+49=", 2) | Out-Null
$d.Content.Find.Execute("88-9=", $true, $true, $false, $false, $false, $true, 1, $false, "59+15=", 2) | Out-Null
$d.Content.Find.Execute("59+25=", $true, $true, $false, $false, $false, $true, 1, $false, "5+75=", 2) | Out-Null
$d.Content.Find.Execute("78-41=", $true, $true, $false, $false, $false, $true, 1, $false, "39+6=", 2) | Out-Null
$d.Content.Find.Execute("47-42=", $true, $true, $false, $false, $false, $true, 1, $false, "36+58=", 2) | Out-Null
$d.Content.Find.Execute("99-53=", $true, $true, $false, $false, $false, $true, 1, $false, "30-12=", 2) | Out-Null
$d.Content.Find.Execute("13+50=", $true, $true, $false, $false, $false, $true, 1, $false, "55-52=", 2) | Out-Null
$d.Content.Find.Execute("97-85=", $true, $true, $false, $false, $false, $true, 1, $false, "91-55=", 2) | Out-Null
$d.Content.Find.Execute("92-42=", $true, $true, $false, $false, $false, $true, 1, $false, "21+61=", 2) | Out-Null
$d.Content.Find.Execute("92-86=", $true, $true, $false, $false, $false, $true, 1, $false, "15+53=", 2) | Out-Null
$d.Content.Find.Execute("72-65=", $true, $true, $false, $false, $false, $true, 1, $false, "80-54=", 2) | Out-Null
$d.Content.Find.Execute("97-80=", $true, $true, $false, $false, $false, $true, 1, $false, "44-13=", 2) | Out-Null
$d.Content.Find.Execute("94-39=", $true, $true, $false, $false, $false, $true, 1, $false, "6+37=", 2) | Out-Null
$d.Content.Find.Execute("53+40=", $true, $true, $false, $false, $false, $true, 1, $false, "24-22=", 2) | Out-Null
$d.Content.Find.Execute("93-53=", $true, $true, $false, $false, $false, $true, 1, $false, "48-14=", 2) | Out-Null
$d.Content.Find.Execute("25+72=", $true, $true, $false, $false, $false, $true, 1, $false, "72-55=", 2) | Out-Null
$d.Content.Find.Execute("10+57=", $true, $true, $false, $false, $false, $true, 1, $false, "94-43=", 2) | Out-Null
$d.Content.Find.Execute("90+2=", $true, $true, $false, $false, $false, $true, 1, $false, "14+25=", 2) | Out-Null
$d.Content.Find.Execute("58-0=", $true, $true, $false, $false, $false, $true, 1, $false, "89-18=", 2) | Out-Null
$d.Content.Find.Execute("94-79=", $true, $true, $false, $false, $false, $true, 1, $false, "99-27=", 2) | Out-Null
$d.Content.Find.Execute("94-25=", $true, $true, $false, $false, $false, $true, 1, $false, "93-25=", 2) | Out-Null
$d.Content.Find.Execute("64-7=", $true, $true, $false, $false, $false, $true, 1, $false, "26+52=", 2) | Out-Null
$d.Content.Find.Execute("21-7=", $true, $true, $false, $false, $false, $true, 1, $false, "23-12=", 2) | Out-Null
$d.Content.Find.Execute("80+18=", $true, $true, $false, $false, $false, $true, 1, $false, "21+19=", 2) | Out-Null
$d.Content.Find.Execute("3+1=", $true, $true, $false, $false, $false, $true, 1, $false, "46-42=", 2) | Out-Null
$d.Content.Find.Execute("83-74=", $true, $true, $false, $false, $false, $true, 1, $false, "33-19=", 2) | Out-Null
$d.Content.Find.Execute("44-12=", $true, $true, $false, $false, $false, $true, 1, $false, "63-32=", 2) | Out-Null
$d.Content.Find.Execute("90+5=", $true, $true, $false, $false, $false, $true, 1, $false, "18+37=", 2) | Out-Null
$d.Content.Find.Execute("88-67=", $true, $true, $false, $false, $false, $true, 1, $false, "84-34=", 2) | Out-Null
$d.Content.Find.Execute("13+16=", $true, $true, $false, $false, $false, $true, 1, $false, "35+15=", 2) | Out-Null
$d.Content.Find.Execute("81-79=", $true, $true, $false, $false, $false, $true, 1, $false, "98-37=", 2) | Out-Null
$d.Content.Find.Execute("44+26=", $true, $true, $false, $false, $false, $true, 1, $false, "53+28=", 2) | Out-Null
$d.Content.Find.Execute("99-71=", $true, $true, $false, $false, $false, $true, 1, $false, "90-36=", 2) | Out-Null
$d.Content.Find.Execute("50+29=", $true, $true, $false, $false, $false, $true, 1, $false, "31+7=", 2) | Out-Null
$d.Content.Find.Execute("61+0=", $true, $true, $false, $false, $false, $true, 1, $false, "78+7=", 2) | Out-Null
$d.Content.Find.Execute("91-83=", $true, $true, $false, $false, $false, $true, 1, $false, "82-20=", 2) | Out-Null
$d.Content.Find.Execute("51+41=", $true, $true, $false, $false, $false, $true, 1, $false, "75-51=", 2) | Out-Null
$d.Content.Find.Execute("35+24=", $true, $true, $false, $false, $false, $true, 1, $false, "76-65=", 2) | Out-Null
$d.Content.Find.Execute("14+8=", $true, $true, $false, $false, $false, $true, 1, $false, "54+5=", 2) | Out-Null
$d.Content.Find.Execute("61+18=", $true, $true, $false, $false, $false, $true, 1, $false, "27-13=", 2) | Out-Null
$d.Content.Find.Execute("7+35=", $true, $true, $false, $false, $false, $true, 1, $false, "6+77=", 2) | Out-Null
$d.Content.Find.Execute("46+12=", $true, $true, $false, $false, $false, $true, 1, $false, "44+22=", 2) | Out-Null
$d.Content.Find.Execute("23-15=", $true, $true, $false, $false, $false, $true, 1, $false, "74-38=", 2) | Out-Null
$d.Content.Find.Execute("23+38=", $true, $true, $false, $false, $false, $true, 1, $false, "34-22=", 2) | Out-Null
$d.Content.Find.Execute("78-59=", $true, $true, $false, $false, $false, $true, 1, $false, "7+59=", 2) | Out-Null
$d.Content.Find.Execute("81+17=", $true, $true, $false, $false, $false, $true, 1, $false, "65-14=", 2) | Out-Null
$d.Content.Find.Execute("75-10=", $true, $true, $false, $false, $false, $true, 1, $false, "61-35=", 2) | Out-Null
$d.Content.Find.Execute("91-5=", $true, $true, $false, $false, $false, $true, 1, $false, "92-7=", 2) | Out-Null
$d.Content.Find.Execute("14+10=", $true, $true, $false, $false, $false, $true, 1, $false, "46-7=", 2) | Out-Null
$d.Content.Find.Execute("79-37=", $true, $true, $false, $false, $false, $true, 1, $false, "17+66=", 2) | Out-Null
$d.Content.Find.Execute("94-22=", $true, $true, $false, $false, $false, $true, 1, $false, "76-67=", 2) | Out-Null
$d.Content.Find.Execute("47+40=", $true, $true, $false, $false, $false, $true, 1, $false, "42+4=", 2) | Out-Null
$d.Content.Find.Execute("17+10=", $true, $true, $false, $false, $false, $true, 1, $false, "61+25=", 2) | Out-Null
$d.Content.Find.Execute("85+8=", $true, $true, $false, $false, $false, $true, 1, $false, "50+34=", 2) | Out-Null
$d.Content.Find.Execute("33+51=", $true, $true, $false, $false, $false, $true, 1, $false, "93-58=", 2) | Out-Null
$d.Content.Find.Execute("49+2=", $true, $true, $false, $false, $false, $true, 1, $false, "33+30=", 2) | Out-Null
$d.Content.Find.Execute("30+46=", $true, $true, $false, $false, $false, $true, 1, $false, "74-40=", 2) | Out-Null
$d.Content.Find.Execute("68+10=", $true, $true, $false, $false, $false, $true, 1, $false, "27+3=", 2) | Out-Null
$d.Content.Find.Execute("29+66=", $true, $true, $false, $false, $false, $true, 1, $false, "14+0=", 2) | Out-Null
$d.Content.Find.Execute("32+10=", $true, $true, $false, $false, $false, $true, 1, $false, "43-33=", 2) | Out-Null
$d.Content.Find.Execute("50-5=", $true, $true, $false, $false, $false, $true, 1, $false, "82-66=", 2) | Out-Null
$d.Content.Find.Execute("61+30=", $true, $true, $false, $false, $false, $true, 1, $false, "65-18=", 2) | Out-Null
$d.Content.Find.Execute("28-1=", $true, $true, $false, $false, $false, $true, 1, $false, "40+42=", 2) | Out-Null
$d.Content.Find.Execute("52-26=", $true, $true, $false, $false, $false, $true, 1, $false, "6+20=", 2) | Out-Null
$d.Content.Find.Execute("36+43=", $true, $true, $false, $false, $false, $true, 1, $false, "48-23=", 2) | Out-Null
$d.Content.Find.Execute("78-68=", $true, $true, $false, $false, $false, $true, 1, $false, "92-45=", 2) | Out-Null
$d.Content.Find.Execute("47-7=", $true, $true, $false, $false, $false, $true, 1, $false, "11+17=", 2) | Out-Null
$d.Content.Find.Execute("10+51=", $true, $true, $false, $false, $false, $true, 1, $false, "1+95=", 2) | Out-Null
$d.Content.Find.Execute("24-10=", $true, $true, $false, $false, $false, $true, 1, $false, "37+21=", 2) | Out-Null
$d.Content.Find.Execute("85-6=", $true, $true, $false, $false, $false, $true, 1, $false, "95-48=", 2) | Out-Null
$d.Content.Find.Execute("28-10=", $true, $true, $false, $false, $false, $true, 1, $false, "61+15=", 2) | Out-Null
$d.Content.Find.Execute("48+14=", $true, $true, $false, $false, $false, $true, 1, $false, "81+4=", 2) | Out-Null
$d.Content.Find.Execute("88-28=", $true, $true, $false, $false, $false, $true, 1, $false, "86-50=", 2) | Out-Null
$d.Content.Find.Execute("75-5=", $true, $true, $false, $false, $false, $true, 1, $false, "98-30=", 2) | Out-Null
$d.Content.Find.Execute("31+28=", $true, $true, $false, $false, $false, $true, 1, $false, "74+17=", 2) | Out-Null
